# Auto-committed on 2023/06/09 週五 16:56:32.10
#
# PfBsDetail.xlsx - "DBD" sheet, rows 21-24 (PerfCnt / PerfAmt / AdjPerfCnt /
# AdjPerfAmt field docs): the remark column (G) is reworked --
#   G21 "追回時為扣除金額後重算之件數" -> "調整後件數"
#   G22 (blank)                         -> "調整後業績金額"   (new remark)
#   G23 "未用(移至房貸專員業績調整檔)"   -> (cleared)
#   G24 "未用(移至房貸專員業績調整檔)"   -> (cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# G21: 追回時為扣除金額後重算之件數 -> 調整後件數
$ws.Range("G21").Value2 = "調整後件數"

# G22 did not have a value before; give it the same look as the rest of the
# column (style copied from G21) before writing its new remark text.
$ws.Range("G21").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value2 = "調整後業績金額"

# G23 / G24: the old "未用(移至房貸專員業績調整檔)" remarks are removed,
# leaving the cells blank (formatting untouched).
$ws.Range("G23").ClearContents()
$ws.Range("G24").ClearContents()

# Reflect the author's final view position/selection on the DBD sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G26").Select()
